$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.000047231523702066599
$ws.Range("C4").Value = 0.0000652277783556201
$ws.Range("D4").Value = 0.98537614727020195
$ws.Range("E4").Value = 0.98554356670379595
$ws.Range("F4").Value = 44.644213973999001
$ws.Range("G4").Value = 43.405451065063403

$ws.Range("B13").Select() | Out-Null
